# "Generate Report for Handoff"
#
# The localization status report marks the e430d54b-...md file as
# "Ready for handoff" (it had been "In Translation"), bumps its
# handoff timestamp(s), and records a new "mt" priority for the
# per-language handoff rows. Three sheets carry this same file's row:
#   - Overview (summary row, one line per source file)
#   - zh-cn     (per-language handoff detail table)
#   - de-de     (per-language handoff detail table)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row for e430d54b-e732-4a17-b978-185b3cc94691.md
# (row 3 - the 6d916e2c file sitting in row 2 is untouched)
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 04:17:21"

# Columns E (zh-cn) and F (de-de) grow to fit the new, longer status text.
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------
# zh-cn sheet: row for e430d54b-e732-4a17-b978-185b3cc94691.md (row 3)
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-06 04:17:16"

# Status column (C) grows to fit the new, longer text.
$zhcn.Columns.Item(3).AutoFit()

# ---------------------------------------------------------------
# de-de sheet: row for e430d54b-e732-4a17-b978-185b3cc94691.md (row 3)
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-06 04:17:21"

# Status column (C) grows to fit the new, longer text.
$dede.Columns.Item(3).AutoFit()
